$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new agenda rows (5-7) for técnico Ryan
$ws.Range("A5").Value = "Ryan"
$ws.Range("B5").Value = "'0673"
$ws.Range("C5").Value = "Antena Rádio Espirita"
$ws.Range("D5").Value = "Sem comunicação de alarmes, câmeras on. Problema no tamper da sirene e parece que o cliente quer ajustar um pouco o ângulo de uma câmera."
$ws.Range("G5").Value = "Pendente"
$ws.Rows.Item(5).RowHeight = 30

$ws.Range("A6").Value = "Ryan"
$ws.Range("B6").Value = "'0645"
$ws.Range("C6").Value = "Obra Rancho Alegre"
$ws.Range("D6").Value = "Local sem comunicação de alarmes já tem uns dias, central comunica só via GPRS. Obra do Paulo."
$ws.Range("G6").Value = "Pendente"
$ws.Rows.Item(6).RowHeight = 30

$ws.Range("A7").Value = "Ryan"
$ws.Range("B7").Value = "'1059"
$ws.Range("C7").Value = "Fundição Carola"
$ws.Range("D7").Value = "Zona 6 aberta, disparando sem parar logo após o arme."
$ws.Range("G7").Value = "Pendente"

# Update the sheet view scroll position and selection (topLeftCell F1, active cell H7)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 6
$ws.Range("H7").Select()
